# Company Work Log update - 12/18/2023 and also use Clockify for worklog
#
# This script reproduces the data edits captured in the target diff:
#  - corrects two "Total Time" numbers (D25, D181)
#  - fills in previously blank Task/Type cells
#  - retypes/relabels a block of Task Description / Type cells for the
#    12/18, 12/19 and 12/20/2023 entries (new task strings such as
#    "NAS Python coding", "use another Worklog app Clockify",
#    "Math simulation on Python", "Raspberry Pi 5V power Layout design",
#    "Raspberry Pi 5V power schemeticupload" get introduced)
#  - leaves the final selection on B222, matching the author's last edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix an earlier day's total time -------------------------------------
$ws.Range("D25").Value = 9

# --- fill in previously empty Task/Type cells -----------------------------
$ws.Range("B36").Value = "LED Ring schematic design"
$ws.Range("C36").Value = "design"

$ws.Range("B78").Value = "LED Ring Layout design"
$ws.Range("C78").Value = "design"

$ws.Range("B99").Value  = "LED Ring Layout design"
$ws.Range("C99").Value  = "design"
$ws.Range("B100").Value = "LED Ring Layout design"
$ws.Range("C100").Value = "design"
$ws.Range("B101").Value = "LED Ring Layout design"
$ws.Range("C101").Value = "design"
$ws.Range("B102").Value = "LED Ring Layout design"
$ws.Range("C102").Value = "design"
$ws.Range("B103").Value = "LED Ring Layout design"
$ws.Range("C103").Value = "design"

$ws.Range("B114").Value = "LED Ring Layout design"
$ws.Range("C114").Value = "design"

$ws.Range("B148").Value = "Raspberry Pi 5V power schemetic design"
$ws.Range("C148").Value = "design"
$ws.Range("B149").Value = "Raspberry Pi 5V power schemetic design"
$ws.Range("C149").Value = "design"
$ws.Range("B150").Value = "Raspberry Pi 5V power schemetic design"
$ws.Range("C150").Value = "design"

# --- retype the 12/13/2023 evening block to "NAS Python coding" ----------
$ws.Range("B171").Value = "NAS Python coding"
$ws.Range("B172").Value = "NAS Python coding"
$ws.Range("B173").Value = "NAS Python coding"
$ws.Range("B174").Value = "NAS Python coding"
$ws.Range("B175").Value = "NAS Python coding"
$ws.Range("B176").Value = "NAS Python coding"

# --- fix the 12/13/2023 total time ----------------------------------------
$ws.Range("D181").Value = 4

$ws.Range("B183").Value = "NAS Python coding"
$ws.Range("C183").Value = "Support"
$ws.Range("B184").Value = "NAS Python coding"
$ws.Range("C184").Value = "Support"
$ws.Range("B185").Value = "NAS Python coding"
$ws.Range("C185").Value = "Support"

# --- 12/18/2023 entries ----------------------------------------------------
$ws.Range("B195").Value = "NAS Python coding"
$ws.Range("C195").Value = "Support"
$ws.Range("B196").Value = "NAS Python coding/UPLOAD"
$ws.Range("C196").Value = "Support"
$ws.Range("B197").Value = "Raspberry Pi 5V power schemetic design"
$ws.Range("B198").Value = "use another Worklog app Clockify"
$ws.Range("C198").Value = "support"
$ws.Range("B199").Value = "use another Worklog app Clockify"
$ws.Range("B200").Value = "use another Worklog app Clockify"
$ws.Range("B201").Value = "Raspberry Pi 5V power schemetic design"
$ws.Range("C201").Value = "design"
$ws.Range("B202").Value = "Raspberry Pi 5V power schemetic design"
$ws.Range("B203").Value = "Raspberry Pi 5V power schemetic design"
$ws.Range("B204").Value = "Raspberry Pi 5V power schemetic design"

# --- 12/19/2023 entries ----------------------------------------------------
$ws.Range("B207").Value = "Raspberry Pi 5V power schemetic design"
$ws.Range("B208").Value = "Raspberry Pi 5V power schemetic design"

# "Math simulation on Python" typed in first for this block, then the
# "Raspberry Pi 5V power Layout design" / "...schemeticupload" rows, which
# keeps the shared-string creation order consistent with the saved file.
$ws.Range("B212").Value = "Math simulation on Python"
$ws.Range("B213").Value = "Math simulation on Python"
$ws.Range("B214").Value = "Math simulation on Python"
$ws.Range("B215").Value = "Math simulation on Python"
$ws.Range("B216").Value = "Math simulation on Python"

$ws.Range("B210").Value = "Raspberry Pi 5V power Layout design"
$ws.Range("B211").Value = "Raspberry Pi 5V power Layout design"

$ws.Range("B209").Value = "Raspberry Pi 5V power schemeticupload"

# --- 12/20/2023 entries ----------------------------------------------------
$ws.Range("B219").Value = "Math simulation on Python"
$ws.Range("C219").Value = "design"
$ws.Range("B220").Value = "Math simulation on Python"
$ws.Range("C220").Value = "design"
$ws.Range("B221").Value = "Math simulation on Python"
$ws.Range("C221").Value = "design"
$ws.Range("B222").Value = "Math simulation on Python"
$ws.Range("C222").Value = "design"

# --- match the author's final view/selection state -------------------------
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 193
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("B222").Select()
